$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Swap the "Periodo Mora" labels between the two data rows.
$ws.Range("E16").Value = "1611"
$ws.Range("E17").Value = "1612"

# Update the "Valor Mora" amounts for both rows.
$ws.Range("G16").Value = 1000000
$ws.Range("G17").Value = 1000000
